# Central_America.xlsx: removed Cayambe because it doesn't work
#
# The "Cayambe" volcano occupies two data rows (direction "D" and "A") in
# Sheet1's table. Deleting those two entire rows shifts every row below
# them up by two, shrinks the worksheet's used range accordingly, and
# drops the now-unused "Cayambe" / its polygon strings from the shared
# string table automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cayambe's two rows are rows 8 and 9 (D then A direction).
$ws.Rows("8:9").Delete()

# Leave the same cell selected/active as in the saved workbook.
$ws.Range("C18").Select() | Out-Null
